$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1806
$ws1.Range("F4").Value = 477
$ws1.Range("F7").Value = 642
$ws1.Range("F9").Value = 1770
$ws1.Range("F10").Value = 385
$ws1.Range("F12").Value = 827
$ws1.Range("F13").Value = 352
$ws1.Range("F15").Value = 12935
$ws1.Range("F16").Value = 12891
$ws1.Range("F18").Value = 751
$ws1.Range("F22").Value = 599
$ws1.Range("F23").Value = 2025
$ws1.Range("F26").Value = 10
$ws1.Range("F28").Value = 111
$ws1.Range("F30").Value = 694

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 85

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 180

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 180
$ws4.Range("F5").Value = 1806
$ws4.Range("F6").Value = 477
$ws4.Range("F11").Value = 642
$ws4.Range("F14").Value = 1770
$ws4.Range("F15").Value = 385
$ws4.Range("F17").Value = 827
$ws4.Range("F18").Value = 352
$ws4.Range("F21").Value = 12935
$ws4.Range("F22").Value = 12891
$ws4.Range("F24").Value = 751
$ws4.Range("F28").Value = 599
$ws4.Range("F31").Value = 2025
$ws4.Range("F34").Value = 10
$ws4.Range("F38").Value = 111
$ws4.Range("F40").Value = 694
$ws4.Range("F41").Value = 85
